# Reconciliação de Dados concluído e Bias corrigido
# Adds the "Reconciliado" column (with recomputed "Polarização (bias)")
# to the Statistics sheet, for both the Tempos (A:G) and Distâncias (J:P)
# blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistics")

# --- Header row (row 1): insert "Reconciliado" column, shifting the
#     existing Polarização (bias) / Precisão / Incerteza labels over ---
$ws.Range("D1").Value = "Reconciliado"
$ws.Range("E1").Value = "Polarização (bias)"
$ws.Range("F1").Value = "Precisão"
$ws.Range("G1").Value = "Incerteza"

$ws.Range("J1").Value = "Distâncias"
$ws.Range("M1").Value = "Reconciliado"
$ws.Range("N1").Value = "Polarização (bias)"
$ws.Range("O1").Value = "Precisão"
$ws.Range("P1").Value = "Incerteza"

$rowData = @{
    2 = @{ D = "18.537839224999992"; E = "0.0"; M = "145.16243614688585"; N = "-2.842170943040401E-14" }
    3 = @{ D = "4.230804285"; E = "0.0"; M = "49.34778874875626"; N = "-1.4210854715202004E-14" }
    4 = @{ D = "28.393963064999987"; E = "-3.552713678800501E-15"; M = "303.98047513760923"; N = "0.0" }
    5 = @{ D = "9.16259248"; E = "-3.552713678800501E-15"; M = "34.944094476626745"; N = "-1.4210854715202004E-14" }
    6 = @{ D = "11.609595815000002"; E = "0.0"; M = "122.00003149564448"; N = "2.842170943040401E-14" }
    7 = @{ D = "8.555156629999995"; E = "-1.7763568394002505E-15"; M = "103.4923986614163"; N = "-2.842170943040401E-14" }
    8 = @{ D = "9.24385705"; E = "0.0"; M = "99.94426092346669"; N = "0.0" }
    9 = @{ D = "1.2541097750000003"; E = "2.220446049250313E-16"; M = "13.67235215682766"; N = "0.0" }
    10 = @{ D = "12.000239910000007"; E = "3.552713678800501E-15"; M = "143.8389105702794"; N = "0.0" }
    11 = @{ D = "2.3623795949999997"; E = "-4.440892098500626E-16"; M = "29.91509970791585"; N = "0.0" }
    12 = @{ D = "22.422399525"; E = "-1.0658141036401503E-14"; M = "178.5787951191696"; N = "-2.842170943040401E-14" }
    13 = @{ D = "1.2286224250000002"; E = "0.0"; M = "7.856928791004396"; N = "8.881784197001252E-16" }
    14 = @{ D = "9.011277795000003"; E = "1.7763568394002505E-15"; M = "106.68960203965672"; N = "-2.842170943040401E-14" }
    15 = @{ D = "12.854472385"; E = "0.0"; M = "146.96217253823554"; N = "2.842170943040401E-14" }
    16 = @{ D = "5.682110414999999"; E = "0.0"; M = "68.40801674510804"; N = "-1.4210854715202004E-14" }
    17 = @{ D = "8.700121075"; E = "0.0"; M = "73.44656897532904"; N = "0.0" }
    18 = @{ D = "11.006027159999997"; E = "-1.7763568394002505E-15"; M = "114.99340137681098"; N = "1.4210854715202004E-14" }
    19 = @{ D = "12.094836520000001"; E = "0.0"; M = "95.59441832196416"; N = "0.0" }
    20 = @{ D = "188.35040512999998"; E = "5.684341886080802E-14"; M = "1838.8277519327066"; N = "-4.547473508864641E-13" }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]

    # Tempos block: column D becomes "Reconciliado" (the reconciled value,
    # equal to the mean), column E becomes the recomputed bias (Reconciliado - Médias).
    $ws.Cells.Item($r, 4).Value = [double]$vals.D
    $ws.Cells.Item($r, 5).Value = [double]$vals.E

    # Distâncias block: column M becomes "Reconciliado", column N the
    # recomputed bias (Reconciliado - Médias).
    $ws.Cells.Item($r, 13).Value = [double]$vals.M
    $ws.Cells.Item($r, 14).Value = [double]$vals.N
}
